# Logged Week 15 and simulated Week 16
$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$rushing = $wb.Worksheets.Item("Rushing")

# Row 2: K.Cousins
$rushing.Cells.Item(2, 3).Value = 3   # C2
$rushing.Cells.Item(2, 5).Value = 3   # E2

# Row 3: D.Cook
$rushing.Cells.Item(3, 3).Value = 148 # C3
$rushing.Cells.Item(3, 4).Value = 97  # D3
$rushing.Cells.Item(3, 5).Value = 22  # E3
$rushing.Cells.Item(3, 6).Value = 43  # F3

# Row 5: K.Nwangwu
$rushing.Cells.Item(5, 3).Value = 2   # C5
$rushing.Cells.Item(5, 4).Value = 2   # D5
$rushing.Cells.Item(5, 6).Value = 1   # F5

# --- Receiving sheet ---
$receiving = $wb.Worksheets.Item("Receiving")

# Row 2: D.Cook
$receiving.Cells.Item(2, 3).Value = 44 # C2
$receiving.Cells.Item(2, 4).Value = 29 # D2

# Row 4: C.Ham
$receiving.Cells.Item(4, 3).Value = 10 # C4
$receiving.Cells.Item(4, 4).Value = 10 # D4

# Row 5: J.Jefferson
$receiving.Cells.Item(5, 3).Value = 97 # C5
$receiving.Cells.Item(5, 4).Value = 68 # D5
$receiving.Cells.Item(5, 5).Value = 45 # E5
$receiving.Cells.Item(5, 7).Value = 19 # G5
$receiving.Cells.Item(5, 8).Value = 8  # H5

# Row 8: K.Osborn
$receiving.Cells.Item(8, 3).Value = 51 # C8
$receiving.Cells.Item(8, 4).Value = 33 # D8

# Row 10: D.Westbrook
$receiving.Cells.Item(10, 3).Value = 63 # C10
$receiving.Cells.Item(10, 4).Value = 46 # D10
$receiving.Cells.Item(10, 5).Value = 6  # E10
$receiving.Cells.Item(10, 7).Value = 14 # G10

# Row 12: I.Smith
$receiving.Cells.Item(12, 3).Value = 2 # C12
$receiving.Cells.Item(12, 4).Value = 2 # D12
$receiving.Cells.Item(12, 5).Value = 1 # E12
$receiving.Cells.Item(12, 7).Value = 1 # G12
$receiving.Cells.Item(12, 8).Value = 1 # H12
